$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SortedEvents")

function Swap-Cells($ws, $ref1, $ref2) {
    $v1 = $ws.Range($ref1).Value()
    $v2 = $ws.Range($ref2).Value()
    $ws.Range($ref1).Value = $v2
    $ws.Range($ref2).Value = $v1
}

# Rows 9 and 10: swap A and B columns
Swap-Cells $ws "A9" "A10"
Swap-Cells $ws "B9" "B10"

# Rows 19 and 20: swap A and B columns
Swap-Cells $ws "A19" "A20"
Swap-Cells $ws "B19" "B20"

# Rows 25 and 26: swap A and B columns
Swap-Cells $ws "A25" "A26"
Swap-Cells $ws "B25" "B26"

# Rows 52 and 53: swap B column only
Swap-Cells $ws "B52" "B53"
